# The workbook's "artfynd" sheet lists 17 species observation records
# (rows 2-18, columns A:AY). The upstream automatic export re-ran and the
# same 17 records came back in a different row order - this is a pure
# permutation of whole data rows; no cell values within a row change and
# the header row (row 1) is untouched.
#
# Mapping: destination row -> source row (both in the ORIGINAL workbook)
#   2 <- 14   3 <- 4    4 <- 2    5 <- 15   6 <- 7    7 <- 11   8 <- 13
#   9 <- 12  10 <- 6   11 <- 10  12 <- 5   13 <- 3   14 <- 8   15 <- 9
#  16 <- 18  17 <- 16  18 <- 17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1    # A
$lastCol  = 51   # AY
$firstRow = 2
$lastRow  = 18

$rowMap = @{
    2  = 14
    3  = 4
    4  = 2
    5  = 15
    6  = 7
    7  = 11
    8  = 13
    9  = 12
    10 = 6
    11 = 10
    12 = 5
    13 = 3
    14 = 8
    15 = 9
    16 = 18
    17 = 16
    18 = 17
}

# 1) Snapshot every source row's values BEFORE writing anything, so that
#    overlapping read/write cycles among the rows can't clobber data that
#    is still needed for a later destination row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += ,($ws.Cells.Item($r, $c).Value())
    }
    $snapshot[$r] = $vals
}

# 2) Write each destination row from the snapshot of its mapped source row.
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - $firstCol]
    }
}
